# Updating test files to match the current format in beta
#
# optimization_parameters sheet:
#   - rename the "Model" label (A8) to "production_function"
#   - insert a new "L_curve" parameter row right after it
#   - drop the obsolete "Deletion" row further down
#   - make this sheet the active / selected tab (it was network_weights before)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Rename "Model" -> "production_function" (value stays "Sigmoid")
$ws.Range("A8").Value = "production_function"

# Insert the new L_curve row directly under production_function/Sigmoid
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").NumberFormat = "0.00E+00"
$ws.Range("B9").Value = 0

# Remove the old "Deletion" row (now shifted down to row 17)
$ws.Rows.Item(17).Delete()

# Make optimization_parameters the active sheet/tab and mirror the saved
# selection (whole row 9, i.e. the new L_curve row)
$ws.Activate()
$ws.Rows.Item(9).Select()
